# Updates cryptos price/volume table (rows 2-51) to the refreshed
# coinranking.com snapshot, matching the Oct 7 2023 GitHub Actions sync.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.924.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.32%  "

# Row 3
$ws.Range("D3").Value = "'1.642.74"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "'213.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.86%  "

# Row 6
$ws.Range("E6").Value = "  -0.17%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "'23.72"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.09%  "

# Row 9
$ws.Range("E9").Value = "  +0.57%  "

# Row 10
$ws.Range("D10").Value = "'0.0617"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.93%  "

# Row 11
$ws.Range("D11").Value = "'0.0875"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.59%  "

# Row 12
$ws.Range("D12").Value = "'1.875.76"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.13%  "

# Row 13
$ws.Range("D13").Value = "'1.650.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.50%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.574"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.45%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'4.07"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.81%  "

# Row 16
$ws.Range("D16").Value = "'65.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.00%  "

# Row 17
$ws.Range("D17").Value = "'27.906.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.41%  "

# Row 18
$ws.Range("D18").Value = "'230.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("E19").Value = "  +0.82%  "

# Row 20
$ws.Range("D20").Value = "'7.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.27%  "

# Row 21
$ws.Range("E21").Value = "  +0.04%  "

# Row 22
$ws.Range("D22").Value = "'11.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.64%  "

# Row 23
$ws.Range("D23").Value = "'4.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.35%  "

# Row 24
$ws.Range("D24").Value = "'2.11"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.24%  "

# Row 25
$ws.Range("D25").Value = "'152.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.43%  "

# Row 26
$ws.Range("E26").Value = "  +0.56%  "

# Row 27
$ws.Range("E27").Value = "  +0.80%  "

# Row 28
$ws.Range("D28").Value = "'15.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.17%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  +1.22%  "

# Row 31
$ws.Range("E31").Value = "  +0.24%  "

# Row 32
$ws.Range("E32").Value = "  +1.89%  "

# Row 33
$ws.Range("D33").Value = "'1.424.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.20%  "

# Row 34
$ws.Range("E34").Value = "  +1.61%  "

# Row 35
$ws.Range("E35").Value = "  +2.09%  "

# Row 36
$ws.Range("D36").Value = "'2.35"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  +1.64%  "

# Row 38
$ws.Range("E38").Value = "  +1.00%  "

# Row 39
$ws.Range("E39").Value = "  -1.36%  "

# Row 40
$ws.Range("D40").Value = "'0.557"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.08%  "

# Row 41
$ws.Range("E41").Value = "  +2.67%  "

# Row 42
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").Value = "'2.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.75%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'66.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.19%  "

# Row 45
$ws.Range("E45").Value = "  +3.30%  "

# Row 46
$ws.Range("D46").Value = "'5.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.19%  "

# Row 47
$ws.Range("E47").Value = "  -0.03%  "

# Row 48
$ws.Range("D48").Value = "'1.784.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.08%  "

# Row 49
$ws.Range("D49").Value = "'88.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.58%  "

# Row 50
$ws.Range("E50").Value = "  +1.03%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0506"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.54%  "
